$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_5_8_23"
$ws.Cells.Item(2, 2).Value = 0.7647329332432279
$ws.Cells.Item(2, 3).Value = 0.7566646430420267
$ws.Cells.Item(2, 4).Value = 0.2480315354003723
$ws.Cells.Item(2, 5).Value = 0.5353912794703288
$ws.Cells.Item(2, 6).Value = 0.2603712975978851
$ws.Cells.Item(2, 7).Value = 0.3901068568229675
$ws.Cells.Item(2, 8).Value = 1.189631700515747
$ws.Cells.Item(2, 9).Value = 0.7663536071777344

$ws.Cells.Item(3, 1).Value = "model_5_8_22"
$ws.Cells.Item(3, 2).Value = 0.7648653049872381
$ws.Cells.Item(3, 3).Value = 0.7564162294686795
$ws.Cells.Item(3, 4).Value = 0.2530414186145908
$ws.Cells.Item(3, 5).Value = 0.537524055422022
$ws.Cells.Item(3, 6).Value = 0.2602248191833496
$ws.Cells.Item(3, 7).Value = 0.3905050754547119
$ws.Cells.Item(3, 8).Value = 1.181705951690674
$ws.Cells.Item(3, 9).Value = 0.762835681438446

$ws.Cells.Item(4, 1).Value = "model_5_8_24"
$ws.Cells.Item(4, 2).Value = 0.7648997684318564
$ws.Cells.Item(4, 3).Value = 0.7566782702451125
$ws.Cells.Item(4, 4).Value = 0.2451880527088666
$ws.Cells.Item(4, 5).Value = 0.5341144665264681
$ws.Cells.Item(4, 6).Value = 0.2601867020130157
$ws.Cells.Item(4, 7).Value = 0.390084981918335
$ws.Cells.Item(4, 8).Value = 1.194130182266235
$ws.Cells.Item(4, 9).Value = 0.768459677696228

$ws.Cells.Item(5, 1).Value = "model_5_8_21"
$ws.Cells.Item(5, 2).Value = 0.7680010020613707
$ws.Cells.Item(5, 3).Value = 0.7582368705345263
$ws.Cells.Item(5, 4).Value = 0.2678336349241986
$ws.Cells.Item(5, 5).Value = 0.5451376049230083
$ws.Cells.Item(5, 6).Value = 0.2567545473575592
$ws.Cells.Item(5, 7).Value = 0.3875862956047058
$ws.Cells.Item(5, 8).Value = 1.158304333686829
$ws.Cells.Item(5, 9).Value = 0.7502774000167847

$ws.Cells.Item(6, 1).Value = "model_5_8_20"
$ws.Cells.Item(6, 2).Value = 0.7681199805223183
$ws.Cells.Item(6, 3).Value = 0.7582318048163502
$ws.Cells.Item(6, 4).Value = 0.2709913962459282
$ws.Cells.Item(6, 5).Value = 0.5465602327873382
$ws.Cells.Item(6, 6).Value = 0.256622850894928
$ws.Cells.Item(6, 7).Value = 0.3875944316387177
$ws.Cells.Item(6, 8).Value = 1.153308629989624
$ws.Cells.Item(6, 9).Value = 0.7479308843612671

$ws.Cells.Item(7, 1).Value = "model_5_8_18"
$ws.Cells.Item(7, 2).Value = 0.7719244605634186
$ws.Cells.Item(7, 3).Value = 0.7588497325191847
$ws.Cells.Item(7, 4).Value = 0.3022836546137829
$ws.Cells.Item(7, 5).Value = 0.5610022997971802
$ws.Cells.Item(7, 6).Value = 0.2524124085903168
$ws.Cells.Item(7, 7).Value = 0.3866037726402283
$ws.Cells.Item(7, 8).Value = 1.103803515434265
$ws.Cells.Item(7, 9).Value = 0.7241092324256897

$ws.Cells.Item(8, 1).Value = "model_5_8_19"
$ws.Cells.Item(8, 2).Value = 0.772156122484139
$ws.Cells.Item(8, 3).Value = 0.7590279455180204
$ws.Cells.Item(8, 4).Value = 0.2985597125523173
$ws.Cells.Item(8, 5).Value = 0.5594130139037725
$ws.Cells.Item(8, 6).Value = 0.2521560192108154
$ws.Cells.Item(8, 7).Value = 0.3863180875778198
$ws.Cells.Item(8, 8).Value = 1.109694957733154
$ws.Cells.Item(8, 9).Value = 0.7267307043075562

$ws.Cells.Item(9, 1).Value = "model_5_8_14"
$ws.Cells.Item(9, 2).Value = 0.7758362653938722
$ws.Cells.Item(9, 3).Value = 0.7531671728775778
$ws.Cells.Item(9, 4).Value = 0.3521614956610498
$ws.Cells.Item(9, 5).Value = 0.5805905253591408
$ws.Cells.Item(9, 6).Value = 0.2480832189321518
$ws.Cells.Item(9, 7).Value = 0.3957138955593109
$ws.Cells.Item(9, 8).Value = 1.024895668029785
$ws.Cells.Item(9, 9).Value = 0.6917992234230042

$ws.Cells.Item(10, 1).Value = "model_5_8_16"
$ws.Cells.Item(10, 2).Value = 0.7760277611904829
$ws.Cells.Item(10, 3).Value = 0.755797648178276
$ws.Cells.Item(10, 4).Value = 0.3388817115335638
$ws.Cells.Item(10, 5).Value = 0.5759504866974523
$ws.Cells.Item(10, 6).Value = 0.2478712797164917
$ws.Cells.Item(10, 7).Value = 0.3914967775344849
$ws.Cells.Item(10, 8).Value = 1.045904636383057
$ws.Cells.Item(10, 9).Value = 0.699452817440033

$ws.Cells.Item(11, 1).Value = "model_5_8_15"
$ws.Cells.Item(11, 2).Value = 0.7760385592692354
$ws.Cells.Item(11, 3).Value = 0.7546941321708422
$ws.Cells.Item(11, 4).Value = 0.3455113185622848
$ws.Cells.Item(11, 5).Value = 0.5783747649427132
$ws.Cells.Item(11, 6).Value = 0.2478593289852142
$ws.Cells.Item(11, 7).Value = 0.3932659029960632
$ws.Cells.Item(11, 8).Value = 1.0354163646698
$ws.Cells.Item(11, 9).Value = 0.6954541206359863

$ws.Cells.Item(12, 1).Value = "model_5_8_17"
$ws.Cells.Item(12, 2).Value = 0.7762914822106416
$ws.Cells.Item(12, 3).Value = 0.7565585813301661
$ws.Cells.Item(12, 4).Value = 0.334675833873798
$ws.Cells.Item(12, 5).Value = 0.5744432854639906
$ws.Cells.Item(12, 6).Value = 0.2475794106721878
$ws.Cells.Item(12, 7).Value = 0.3902768790721893
$ws.Cells.Item(12, 8).Value = 1.052558422088623
$ws.Cells.Item(12, 9).Value = 0.7019389271736145

$ws.Cells.Item(13, 1).Value = "model_5_8_0"
$ws.Cells.Item(13, 2).Value = 0.7807764127470322
$ws.Cells.Item(13, 3).Value = 0.6934840448288333
$ws.Cells.Item(13, 4).Value = 0.5672646570924258
$ws.Cells.Item(13, 5).Value = 0.6469668427976798
$ws.Cells.Item(13, 6).Value = 0.2426159232854843
$ws.Cells.Item(13, 7).Value = 0.4913958311080933
$ws.Cells.Item(13, 8).Value = 0.6845974326133728
$ws.Cells.Item(13, 9).Value = 0.5823140740394592

$ws.Cells.Item(14, 1).Value = "model_5_8_1"
$ws.Cells.Item(14, 2).Value = 0.7813261138547781
$ws.Cells.Item(14, 3).Value = 0.6789378661338734
$ws.Cells.Item(14, 4).Value = 0.5670294511732765
$ws.Cells.Item(14, 5).Value = 0.6393760868508569
$ws.Cells.Item(14, 6).Value = 0.2420075684785843
$ws.Cells.Item(14, 7).Value = 0.5147157311439514
$ws.Cells.Item(14, 8).Value = 0.6849695444107056
$ws.Cells.Item(14, 9).Value = 0.5948348045349121

$ws.Cells.Item(15, 1).Value = "model_5_8_13"
$ws.Cells.Item(15, 2).Value = 0.7816654233365891
$ws.Cells.Item(15, 3).Value = 0.7442782025610442
$ws.Cells.Item(15, 4).Value = 0.4010557431471978
$ws.Cells.Item(15, 5).Value = 0.5980847463230582
$ws.Cells.Item(15, 6).Value = 0.2416320592164993
$ws.Cells.Item(15, 7).Value = 0.409964382648468
$ws.Cells.Item(15, 8).Value = 0.9475437998771667
$ws.Cells.Item(15, 9).Value = 0.6629432439804077

$ws.Cells.Item(16, 1).Value = "model_5_8_12"
$ws.Cells.Item(16, 2).Value = 0.7818410857513254
$ws.Cells.Item(16, 3).Value = 0.7417840534979725
$ws.Cells.Item(16, 4).Value = 0.4084717952419112
$ws.Cells.Item(16, 5).Value = 0.6001489238254268
$ws.Cells.Item(16, 6).Value = 0.2414376437664032
$ws.Cells.Item(16, 7).Value = 0.4139629006385803
$ws.Cells.Item(16, 8).Value = 0.9358114004135132
$ws.Cells.Item(16, 9).Value = 0.6595384478569031

$ws.Cells.Item(17, 1).Value = "model_5_8_2"
$ws.Cells.Item(17, 2).Value = 0.7832452998931489
$ws.Cells.Item(17, 3).Value = 0.6794030229032146
$ws.Cells.Item(17, 4).Value = 0.5656115380652724
$ws.Cells.Item(17, 5).Value = 0.6389757465364482
$ws.Cells.Item(17, 6).Value = 0.2398835718631744
$ws.Cells.Item(17, 7).Value = 0.5139700174331665
$ws.Cells.Item(17, 8).Value = 0.6872127056121826
$ws.Cells.Item(17, 9).Value = 0.5954951047897339

$ws.Cells.Item(18, 1).Value = "model_5_8_11"
$ws.Cells.Item(18, 2).Value = 0.7833395599614169
$ws.Cells.Item(18, 3).Value = 0.7369423634091525
$ws.Cells.Item(18, 4).Value = 0.4260915482700741
$ws.Cells.Item(18, 5).Value = 0.605609772589536
$ws.Cells.Item(18, 6).Value = 0.2397792637348175
$ws.Cells.Item(18, 7).Value = 0.421724945306778
$ws.Cells.Item(18, 8).Value = 0.9079365730285645
$ws.Cells.Item(18, 9).Value = 0.650530993938446

$ws.Cells.Item(19, 1).Value = "model_5_8_10"
$ws.Cells.Item(19, 2).Value = 0.7833507425057902
$ws.Cells.Item(19, 3).Value = 0.731592204387874
$ws.Cells.Item(19, 4).Value = 0.4388556006778394
$ws.Cells.Item(19, 5).Value = 0.6086184449836627
$ws.Cells.Item(19, 6).Value = 0.2397668659687042
$ws.Cells.Item(19, 7).Value = 0.4303021132946014
$ws.Cells.Item(19, 8).Value = 0.8877435326576233
$ws.Cells.Item(19, 9).Value = 0.645568311214447

$ws.Cells.Item(20, 1).Value = "model_5_8_9"
$ws.Cells.Item(20, 2).Value = 0.7838845596768413
$ws.Cells.Item(20, 3).Value = 0.7262023740038513
$ws.Cells.Item(20, 4).Value = 0.4546830417148584
$ws.Cells.Item(20, 5).Value = 0.6129886454819
$ws.Cells.Item(20, 6).Value = 0.2391761094331741
$ws.Cells.Item(20, 7).Value = 0.4389429092407227
$ws.Cells.Item(20, 8).Value = 0.8627041578292847
$ws.Cells.Item(20, 9).Value = 0.6383597850799561

$ws.Cells.Item(21, 1).Value = "model_5_8_3"
$ws.Cells.Item(21, 2).Value = 0.7851205578612916
$ws.Cells.Item(21, 3).Value = 0.6815617222200214
$ws.Cells.Item(21, 4).Value = 0.5622346762865544
$ws.Cells.Item(21, 5).Value = 0.6385619398891964
$ws.Cells.Item(21, 6).Value = 0.2378082275390625
$ws.Cells.Item(21, 7).Value = 0.5105092525482178
$ws.Cells.Item(21, 8).Value = 0.6925549507141113
$ws.Cells.Item(21, 9).Value = 0.5961776971817017

$ws.Cells.Item(22, 1).Value = "model_5_8_4"
$ws.Cells.Item(22, 2).Value = 0.7863636937177041
$ws.Cells.Item(22, 3).Value = 0.6839210510187617
$ws.Cells.Item(22, 4).Value = 0.5564299505243746
$ws.Cells.Item(22, 5).Value = 0.6371557656618937
$ws.Cells.Item(22, 6).Value = 0.2364324331283569
$ws.Cells.Item(22, 7).Value = 0.506726861000061
$ws.Cells.Item(22, 8).Value = 0.7017382383346558
$ws.Cells.Item(22, 9).Value = 0.5984970927238464

$ws.Cells.Item(23, 1).Value = "model_5_8_5"
$ws.Cells.Item(23, 2).Value = 0.7876578414763966
$ws.Cells.Item(23, 3).Value = 0.695114888927473
$ws.Cells.Item(23, 4).Value = 0.5485534367625118
$ws.Cells.Item(23, 5).Value = 0.6393608187656739
$ws.Cells.Item(23, 6).Value = 0.2350002229213715
$ws.Cells.Item(23, 7).Value = 0.4887812435626984
$ws.Cells.Item(23, 8).Value = 0.7141990661621094
$ws.Cells.Item(23, 9).Value = 0.5948599576950073

$ws.Cells.Item(24, 1).Value = "model_5_8_8"
$ws.Cells.Item(24, 2).Value = 0.7880780574699989
$ws.Cells.Item(24, 3).Value = 0.7067481519802066
$ws.Cells.Item(24, 4).Value = 0.5172043494888989
$ws.Cells.Item(24, 5).Value = 0.6311968854051007
$ws.Cells.Item(24, 6).Value = 0.2345351427793503
$ws.Cells.Item(24, 7).Value = 0.4701312780380249
$ws.Cells.Item(24, 8).Value = 0.7637940049171448
$ws.Cells.Item(24, 9).Value = 0.6083260178565979

$ws.Cells.Item(25, 1).Value = "model_5_8_7"
$ws.Cells.Item(25, 2).Value = 0.7884186630572372
$ws.Cells.Item(25, 3).Value = 0.7017711879389846
$ws.Cells.Item(25, 4).Value = 0.5297793737653704
$ws.Cells.Item(25, 5).Value = 0.6343118581732018
$ws.Cells.Item(25, 6).Value = 0.2341581881046295
$ws.Cells.Item(25, 7).Value = 0.4781101644039154
$ws.Cells.Item(25, 8).Value = 0.7439000606536865
$ws.Cells.Item(25, 9).Value = 0.6031880378723145

$ws.Cells.Item(26, 1).Value = "model_5_8_6"
$ws.Cells.Item(26, 2).Value = 0.7890088516354561
$ws.Cells.Item(26, 3).Value = 0.6979147670129375
$ws.Cells.Item(26, 4).Value = 0.5430661595602401
$ws.Cells.Item(26, 5).Value = 0.6383246692422815
$ws.Cells.Item(26, 6).Value = 0.2335050255060196
$ws.Cells.Item(26, 7).Value = 0.4842926561832428
$ws.Cells.Item(26, 8).Value = 0.7228800654411316
$ws.Cells.Item(26, 9).Value = 0.5965690612792969
